$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Journal")

# New diary entries - row 30
$ws.Range("A30").Value2 = 45523
$ws.Range("B30").Value2 = 0.4375
$ws.Range("C30").Value2 = 0.66666666666666663
$ws.Range("E30").Value2 = "Réalisation de l'application "

# New diary entries - row 31
$ws.Range("A31").Value2 = 45527
$ws.Range("B31").Value2 = 0.35416666666666669
$ws.Range("C31").Value2 = 0.59097222222222223
$ws.Range("E31").Value2 = "Réalisation de l'application "

# Update the selection / scroll position to reflect where the user ended up
$ws.Range("H8").Select()
